$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$labels = @("1M","2M","3M","4M","5M","6M","7M","8M","9M","10M","11M","1Y","15M","18M","21M","2Y","3Y","4Y","5Y","6Y","7Y","8Y","9Y","10Y","11Y","12Y","15Y","20Y","25Y","30Y","40Y","50Y","60Y")

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $labels[$i]
}

$ws.Range("A35").Select()
